$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 345 (shifts existing rows 345-381 down to 347-383)
$ws.Range("A345:A346").EntireRow.Insert()

$newDate = Get-Date -Year 2021 -Month 9 -Day 10 -Hour 0 -Minute 0 -Second 0

# Row 345: new "1a (guarda)" record
$ws.Cells.Item(345, 1).Value2 = 8
$ws.Cells.Item(345, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(345, 3).Value2 = "Coquimbo"
$ws.Cells.Item(345, 4).Value2 = $newDate
$ws.Cells.Item(345, 5).Value2 = 4
$ws.Cells.Item(345, 6).Value2 = 100112004
$ws.Cells.Item(345, 7).Value2 = "Cebolla"
$ws.Cells.Item(345, 8).Value2 = "Sin especificar"
$ws.Cells.Item(345, 9).Value2 = "1a (guarda)"
$ws.Cells.Item(345, 10).Value2 = 3400
$ws.Cells.Item(345, 11).Value2 = 5300
$ws.Cells.Item(345, 12).Value2 = 5500
$ws.Cells.Item(345, 13).Value2 = 5400
$ws.Cells.Item(345, 14).Value2 = "`$/malla 16 kilos"
$ws.Cells.Item(345, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(345, 16).Value2 = 338
$ws.Cells.Item(345, 17).Value2 = 16
$ws.Cells.Item(345, 18).Value2 = "Hortaliza"

# Row 346: new "2a (guarda)" record
$ws.Cells.Item(346, 1).Value2 = 8
$ws.Cells.Item(346, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(346, 3).Value2 = "Coquimbo"
$ws.Cells.Item(346, 4).Value2 = $newDate
$ws.Cells.Item(346, 5).Value2 = 4
$ws.Cells.Item(346, 6).Value2 = 100112004
$ws.Cells.Item(346, 7).Value2 = "Cebolla"
$ws.Cells.Item(346, 8).Value2 = "Sin especificar"
$ws.Cells.Item(346, 9).Value2 = "2a (guarda)"
$ws.Cells.Item(346, 10).Value2 = 1800
$ws.Cells.Item(346, 11).Value2 = 4800
$ws.Cells.Item(346, 12).Value2 = 5000
$ws.Cells.Item(346, 13).Value2 = 4900
$ws.Cells.Item(346, 14).Value2 = "`$/malla 16 kilos"
$ws.Cells.Item(346, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(346, 16).Value2 = 306
$ws.Cells.Item(346, 17).Value2 = 16
$ws.Cells.Item(346, 18).Value2 = "Hortaliza"
